$wb = $excel.ActiveWorkbook

# Sheet "Hoja1": update the conversion text in cell A1
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.37 = 29977.88 pesos`n✅ 29977.88 pesos = 7.37 = 957.76 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Sheet "tasas": update rate cells
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 135.6
$wsTasas.Range("O10").Value = 4065
$wsTasas.Range("N12").Value = 4069
$wsTasas.Range("O12").Value = 130
